$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: replace the run-content of a paragraph (identified by its 1-based
# Paragraphs index at the time of the call) with a fresh set of <w:r> runs,
# leaving the paragraph's own <w:pPr> (style/formatting) and paragraph mark
# untouched.
# ---------------------------------------------------------------------------
function Set-ParaRuns($paraIndex, $innerRunsXml) {
    $p = $d.Paragraphs($paraIndex)
    $full = $p.Range
    $start = $full.Start
    $end = $full.End - 1
    $r = $d.Range($start, $end)
    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $innerRunsXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $r.InsertXML($xml)
}

# 1) "multi do" heading: merge "multi "+"do"+" " into one run, and
#    "<"+"command"+">" into one run.
Set-ParaRuns 18 '<w:r><w:t xml:space="preserve">multi do </w:t></w:r><w:r><w:t>"</w:t></w:r><w:r><w:t>&lt;command&gt;</w:t></w:r><w:r><w:t>"</w:t></w:r>'

# 2) "multi graph" heading: merge "multi "+"graph" into one run.
Set-ParaRuns 23 '<w:r><w:t>multi graph</w:t></w:r>'

# 3) "multi merge" heading: merge "multi "+"merge"+" <ref> [--latest|--exact]"
#    into one run.
Set-ParaRuns 29 '<w:r><w:t>multi merge &lt;ref&gt; [--latest|--exact]</w:t></w:r>'

# 4) "multi open" heading: merge all the separate runs into one run.
Set-ParaRuns 33 '<w:r><w:t>multi open [--all (default)|--main|--deps]</w:t></w:r>'

# 5) "multi update" heading: split into 4 runs adding the new flags.
Set-ParaRuns 41 '<w:r><w:t xml:space="preserve">multi update </w:t></w:r><w:r><w:t>[--all (default)|--main|--deps]</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>[--force] [--commit]</w:t></w:r>'

# ---------------------------------------------------------------------------
# 6) Insert three new "Option" style paragraphs (--all / --main / --deps)
#    right after the ("Removes the specified dependency from multirepo.")
#    paragraph that follows the "multi update" heading (paragraph 42), and
#    before the "--force" option paragraph. The last new paragraph also
#    carries the relocated _GoBack bookmark.
# ---------------------------------------------------------------------------
$insertAfter = $d.Paragraphs(42)
# Insert just *inside* the end of the "Removes..." paragraph (i.e. immediately
# before its paragraph mark) rather than exactly at the boundary between it
# and the next paragraph - InsertXML at an exact paragraph boundary merges
# the final inserted paragraph into whatever followed.
$insertPoint = $d.Range($insertAfter.Range.End - 1, $insertAfter.Range.End - 1)

$newParasXml = @'
<w:p><w:pPr><w:pStyle w:val="Option"/></w:pPr><w:r><w:rPr><w:b/><w:color w:val="E76618" w:themeColor="accent4"/></w:rPr><w:t>--</w:t></w:r><w:r><w:rPr><w:b/><w:color w:val="E76618" w:themeColor="accent4"/></w:rPr><w:t>all</w:t></w:r><w:r><w:tab/></w:r><w:r><w:t xml:space="preserve">Update </w:t></w:r><w:r><w:t>the main repository and all dependencies.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Option"/></w:pPr><w:r><w:rPr><w:b/><w:color w:val="E76618" w:themeColor="accent4"/></w:rPr><w:t>--</w:t></w:r><w:r><w:rPr><w:b/><w:color w:val="E76618" w:themeColor="accent4"/></w:rPr><w:t>main</w:t></w:r><w:r><w:tab/></w:r><w:r><w:t xml:space="preserve">Update </w:t></w:r><w:r><w:t>the main repository.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Option"/></w:pPr><w:r><w:rPr><w:b/><w:color w:val="E76618" w:themeColor="accent4"/></w:rPr><w:t>--</w:t></w:r><w:r><w:rPr><w:b/><w:color w:val="E76618" w:themeColor="accent4"/></w:rPr><w:t>deps</w:t></w:r><w:r><w:tab/></w:r><w:r><w:t>Update</w:t></w:r><w:r><w:t xml:space="preserve"> dependencies</w:t></w:r><w:r><w:t>.</w:t></w:r><w:bookmarkStart w:id="100" w:name="_GoBack_NEW"/><w:bookmarkEnd w:id="100"/></w:p>
'@

$xmlDoc = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + '<w:body>' + $newParasXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertPoint.InsertXML($xmlDoc)

# ---------------------------------------------------------------------------
# 7) Move the _GoBack bookmark: remove it from the early empty Heading1
#    paragraph and place it (already inserted above, under a temp name) -
#    rename the temp bookmark to _GoBack after deleting the original.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

if ($d.Bookmarks.Exists("_GoBack_NEW")) {
    $tmp = $d.Bookmarks("_GoBack_NEW")
    $tmpRange = $tmp.Range
    $tmp.Delete()
    $d.Bookmarks.Add("_GoBack", $tmpRange) | Out-Null
}

Write-Output "done"
